$d = $word.ActiveDocument
$replacements = @(
    @("[[PERSON_1]] – „s [[PERSON_2]]“, „o [[PERSON_3]]“", "[[PERSON_1]] – „s [[PERSON_2]]“, „o [[PERSON_2]]“"),
    @("[[PERSON_4]] – „k Evě Marečkové“, „u [[PERSON_4]]“", "[[PERSON_3]] – „k Evě Marečkové“, „u [[PERSON_3]]“"),
    @("[[PERSON_5]] – „pro [[PERSON_5]]“, „s [[PERSON_5]]“", "[[PERSON_4]] – „pro [[PERSON_4]]“, „s [[PERSON_4]]“"),
    @("[[PERSON_6]] – „s [[PERSON_6]]“, „o [[PERSON_6]]“", "[[PERSON_5]] – „s [[PERSON_5]]“, „o [[PERSON_5]]“"),
    @("[[PERSON_7]] – „u [[PERSON_7]]“, „k [[PERSON_7]]“", "[[PERSON_6]] – „u [[PERSON_6]]“, „k [[PERSON_6]]“"),
    @("[[PERSON_8]] – „o [[PERSON_9]]“, „se [[PERSON_8]]“", "[[PERSON_7]] – „o [[PERSON_7]]“, „se [[PERSON_7]]“"),
    @("[[PERSON_10]] – „k [[PERSON_11]]“, „u [[PERSON_10]]“", "[[PERSON_8]] – „k [[PERSON_9]]“, „u [[PERSON_8]]“"),
    @("[[PERSON_12]] – „s [[PERSON_12]]“, „o [[PERSON_12]]“", "[[PERSON_10]] – „s [[PERSON_10]]“, „o [[PERSON_10]]“"),
    @("[[PERSON_13]] – „u [[PERSON_13]]“, „s [[PERSON_13]]“", "[[PERSON_11]] – „u [[PERSON_11]]“, „s [[PERSON_11]]“"),
    @("[[PERSON_14]] – „s [[PERSON_14]]“, „k [[PERSON_15]]“", "[[PERSON_12]] – „s [[PERSON_12]]“, „k [[PERSON_12]]“"),
    @("[[PERSON_16]] – „s [[PERSON_17]]“, „o [[PERSON_18]]“", "[[PERSON_13]] – „s [[PERSON_14]]“, „o [[PERSON_15]]“"),
    @("[[PERSON_19]] – „ke [[PERSON_20]]“, „o [[PERSON_20]]“", "[[PERSON_16]] – „ke [[PERSON_16]]“, „o [[PERSON_16]]“"),
    @("[[PERSON_21]] – „o [[PERSON_22]]“, „s [[PERSON_21]]“", "[[PERSON_17]] – „o [[PERSON_18]]“, „s [[PERSON_17]]“"),
    @("[[PERSON_23]] – „u [[PERSON_23]]“, „s [[PERSON_23]]“", "[[PERSON_19]] – „u [[PERSON_19]]“, „s [[PERSON_19]]“"),
    @("[[PERSON_24]] – „ke [[PERSON_24]]“, „o [[PERSON_24]]“", "[[PERSON_20]] – „ke [[PERSON_20]]“, „o [[PERSON_20]]“"),
    @("[[PERSON_25]] – „s [[PERSON_25]]“, „o [[PERSON_25]]“", "[[PERSON_21]] – „s [[PERSON_21]]“, „o [[PERSON_21]]“"),
    @("[[PERSON_26]] – „s [[PERSON_26]]“, „o [[PERSON_27]]“", "[[PERSON_22]] – „s [[PERSON_22]]“, „o [[PERSON_22]]“"),
    @("[[PERSON_28]] – „k [[PERSON_29]]“, „od [[PERSON_28]]“", "[[PERSON_23]] – „k [[PERSON_23]]“, „od [[PERSON_23]]“"),
    @("[[PERSON_30]] – „o [[PERSON_31]]“, „s [[PERSON_30]]“", "[[PERSON_24]] – „o [[PERSON_24]]“, „s [[PERSON_24]]“"),
    @("[[PERSON_32]] – „o [[PERSON_33]]“, „se [[PERSON_34]]“", "[[PERSON_25]] – „o [[PERSON_26]]“, „se [[PERSON_27]]“"),
    @("[[PERSON_35]] – „s [[PERSON_35]]“, „u [[PERSON_35]]“", "[[PERSON_28]] – „s [[PERSON_28]]“, „u [[PERSON_28]]“"),
    @("[[PERSON_36]] – „o [[PERSON_37]]“, „s [[PERSON_38]]“", "[[PERSON_29]] – „o [[PERSON_30]]“, „s [[PERSON_31]]“"),
    @("[[PERSON_39]] – „k [[PERSON_39]]“, „o [[PERSON_39]]“", "[[PERSON_32]] – „k [[PERSON_32]]“, „o [[PERSON_32]]“"),
    @("[[PERSON_40]] – „se [[PERSON_40]]“, „o Soně Mikulkové“", "[[PERSON_33]] – „se [[PERSON_33]]“, „o Soně Mikulkové“"),
    @("[[PERSON_41]] – „o [[PERSON_41]]“, „s [[PERSON_41]]“", "[[PERSON_34]] – „o [[PERSON_34]]“, „s [[PERSON_34]]“"),
    @("[[PERSON_42]] – „s [[PERSON_43]]“, „o [[PERSON_42]]“", "[[PERSON_35]] – „s [[PERSON_35]]“, „o [[PERSON_35]]“"),
    @("[[PERSON_44]] – „k [[PERSON_45]]“, „s [[PERSON_44]]“", "[[PERSON_36]] – „k [[PERSON_37]]“, „s [[PERSON_36]]“"),
    @("[[PERSON_46]] – „s [[PERSON_46]]“, „o [[PERSON_47]]“", "[[PERSON_38]] – „s [[PERSON_38]]“, „o [[PERSON_39]]“"),
    @("[[PERSON_48]] – „od [[PERSON_48]]“, „s [[PERSON_49]]“", "[[PERSON_40]] – „od [[PERSON_40]]“, „s [[PERSON_41]]“"),
    @("[[PERSON_50]] – „k [[PERSON_51]]“, „o [[PERSON_51]]“", "[[PERSON_42]] – „k [[PERSON_43]]“, „o [[PERSON_43]]“"),
    @("[[PERSON_52]] – „o [[PERSON_53]]“, „s [[PERSON_54]]“", "[[PERSON_44]] – „o [[PERSON_45]]“, „s [[PERSON_44]]“"),
    @("[[PERSON_55]] – „s [[PERSON_56]]“, „o [[PERSON_57]]“", "[[PERSON_46]] – „s [[PERSON_46]]“, „o [[PERSON_47]]“"),
    @("[[PERSON_58]] – „s [[PERSON_58]]“, „o [[PERSON_59]]“", "[[PERSON_48]] – „s [[PERSON_48]]“, „o [[PERSON_49]]“"),
    @("[[PERSON_60]] – „k [[PERSON_60]]“, „s [[PERSON_61]]“", "[[PERSON_50]] – „k [[PERSON_50]]“, „s [[PERSON_50]]“"),
    @("[[PERSON_62]] – „pro [[PERSON_63]]“, „o [[PERSON_64]]“", "[[PERSON_51]] – „pro [[PERSON_52]]“, „o [[PERSON_53]]“"),
    @("[[PERSON_65]] – „k [[PERSON_65]]“, „o [[PERSON_65]]“", "[[PERSON_54]] – „k [[PERSON_54]]“, „o [[PERSON_54]]“"),
    @("[[PERSON_66]] – „o [[PERSON_67]]“, „s [[PERSON_66]]“", "[[PERSON_55]] – „o [[PERSON_56]]“, „s [[PERSON_55]]“"),
    @("[[PERSON_68]] – „s [[PERSON_69]]“, „o [[PERSON_70]]“", "[[PERSON_57]] – „s [[PERSON_57]]“, „o [[PERSON_58]]“"),
    @("[[PERSON_71]] – „s [[PERSON_71]]“, „o [[PERSON_71]]“", "[[PERSON_59]] – „s [[PERSON_59]]“, „o [[PERSON_59]]“"),
    @("[[PERSON_72]] – „u [[PERSON_72]]“, „o [[PERSON_73]]“", "[[PERSON_60]] – „u [[PERSON_60]]“, „o [[PERSON_61]]“"),
    @("[[PERSON_74]] – „se [[PERSON_75]]“, „o [[PERSON_74]]“", "[[PERSON_62]] – „se [[PERSON_62]]“, „o [[PERSON_62]]“"),
    @("[[PERSON_76]] – „o [[PERSON_77]]“, „s [[PERSON_78]]“", "[[PERSON_63]] – „o [[PERSON_64]]“, „s [[PERSON_65]]“"),
    @("[[PERSON_79]] – „k [[PERSON_80]]“, „o [[PERSON_80]]“", "[[PERSON_66]] – „k [[PERSON_67]]“, „o [[PERSON_67]]“"),
    @("[[PERSON_81]] – „o [[PERSON_82]]“, „s [[PERSON_81]]“", "[[PERSON_68]] – „o [[PERSON_69]]“, „s [[PERSON_68]]“"),
    @("[[PERSON_83]] – „s [[PERSON_83]]“, „o [[PERSON_83]]“", "[[PERSON_70]] – „s [[PERSON_70]]“, „o [[PERSON_70]]“"),
    @("[[PERSON_84]] – „s [[PERSON_84]]“, „o [[PERSON_85]]“", "[[PERSON_71]] – „s [[PERSON_71]]“, „o [[PERSON_72]]“"),
    @("[[PERSON_86]] – „o [[PERSON_87]]“, „s [[PERSON_88]]“", "[[PERSON_73]] – „o [[PERSON_74]]“, „s [[PERSON_74]]“"),
    @("[[PERSON_89]] – „s [[PERSON_90]]“, „o [[PERSON_91]]“", "[[PERSON_75]] – „s [[PERSON_76]]“, „o [[PERSON_77]]“"),
    @("[[PERSON_92]] – „o [[PERSON_93]]“, „s [[PERSON_92]]“", "[[PERSON_78]] – „o [[PERSON_79]]“, „s [[PERSON_78]]“"),
    @("[[PERSON_94]] – „s [[PERSON_95]]“, „o [[PERSON_96]]“", "[[PERSON_80]] – „s [[PERSON_81]]“, „o [[PERSON_82]]“"),
)

$notFound = 0
foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        $notFound = $notFound + 1
        Write-Host "NOT FOUND: $old"
    }
}
Write-Host "Done. Not found count: $notFound"
